$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Test2"
$ws.Range("B3").Value = "A"
$ws.Range("C3").Value = "B"
$ws.Range("D3").Value = "C"
$ws.Range("E3").Value = "D"
$ws.Range("F3").Value = "B"

$ws.Range("F3").Select()
